{"js": "// Update the date header and the 25 two-digit-by-two-digit multiplication\n// problems in the practice sheet. Each old value is unique in the document,\n// so a straightforward search-and-replace (in document order) is safe.\nconst replacements = [\n  [\"2025-02-10 Monday\", \"2025-02-11 Tuesday\"],\n  [\"27\u00d781=\", \"46\u00d747=\"],\n  [\"35\u00d791=\", \"46\u00d791=\"],\n  [\"50\u00d717=\", \"35\u00d712=\"],\n  [\"21\u00d755=\", \"23\u00d727=\"],\n  [\"99\u00d715=\", \"19\u00d798=\"],\n  [\"15\u00d782=\", \"57\u00d764=\"],\n  [\"36\u00d739=\", \"23\u00d765=\"],\n  [\"73\u00d766=\", \"73\u00d737=\"],\n  [\"48\u00d756=\", \"49\u00d741=\"],\n  [\"65\u00d735=\", \"92\u00d730=\"],\n  [\"48\u00d751=\", \"31\u00d740=\"],\n  [\"27\u00d769=\", \"26\u00d752=\"],\n  [\"83\u00d718=\", \"74\u00d723=\"],\n  [\"71\u00d729=\", \"73\u00d788=\"],\n  [\"16\u00d722=\", \"95\u00d733=\"],\n  [\"84\u00d799=\", \"37\u00d711=\"],\n  [\"23\u00d792=\", \"88\u00d791=\"],\n  [\"63\u00d724=\", \"23\u00d730=\"],\n  [\"20\u00d759=\", \"31\u00d783=\"],\n  [\"70\u00d753=\", \"29\u00d787=\"],\n  [\"86\u00d782=\", \"98\u00d798=\"],\n  [\"97\u00d715=\", \"46\u00d780=\"],\n  [\"14\u00d768=\", \"63\u00d778=\"],\n  [\"45\u00d797=\", \"21\u00d774=\"],\n  [\"41\u00d765=\", \"15\u00d728=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date header and the 25 two-digit-by-two-digit multiplication\n# problems in the practice sheet. Each old value is unique in the document,\n# so straightforward Find/Replace (wdReplaceOne) per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-10 Monday\", \"2025-02-11 Tuesday\"),\n    @(\"27\u00d781=\", \"46\u00d747=\"),\n    @(\"35\u00d791=\", \"46\u00d791=\"),\n    @(\"50\u00d717=\", \"35\u00d712=\"),\n    @(\"21\u00d755=\", \"23\u00d727=\"),\n    @(\"99\u00d715=\", \"19\u00d798=\"),\n    @(\"15\u00d782=\", \"57\u00d764=\"),\n    @(\"36\u00d739=\", \"23\u00d765=\"),\n    @(\"73\u00d766=\", \"73\u00d737=\"),\n    @(\"48\u00d756=\", \"49\u00d741=\"),\n    @(\"65\u00d735=\", \"92\u00d730=\"),\n    @(\"48\u00d751=\", \"31\u00d740=\"),\n    @(\"27\u00d769=\", \"26\u00d752=\"),\n    @(\"83\u00d718=\", \"74\u00d723=\"),\n    @(\"71\u00d729=\", \"73\u00d788=\"),\n    @(\"16\u00d722=\", \"95\u00d733=\"),\n    @(\"84\u00d799=\", \"37\u00d711=\"),\n    @(\"23\u00d792=\", \"88\u00d791=\"),\n    @(\"63\u00d724=\", \"23\u00d730=\"),\n    @(\"20\u00d759=\", \"31\u00d783=\"),\n    @(\"70\u00d753=\", \"29\u00d787=\"),\n    @(\"86\u00d782=\", \"98\u00d798=\"),\n    @(\"97\u00d715=\", \"46\u00d780=\"),\n    @(\"14\u00d768=\", \"63\u00d778=\"),\n    @(\"45\u00d797=\", \"21\u00d774=\"),\n    @(\"41\u00d765=\", \"15\u00d728=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
